$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the "Disaster" data column
$ws.Range("I1").Value = "Disaster"

# Values for the new Disaster column (stored as plain decimal numbers)
$ws.Range("I2").Value = 0.025
$ws.Range("I3").Value = 0.05
$ws.Range("I4").Value = 0.1
$ws.Range("I5").Value = 0.075
$ws.Range("I6").Value = 0.15
$ws.Range("I7").Value = 0.13
$ws.Range("I7").Style = "Normal"

# Update the active selection to I8
$ws.Range("I8").Select()
